$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update grade cells in row 14 ------------------------------------------
# P14: "3+" -> "/"
$ws.Range("P14").Value = "/"
# Q14: "3-" -> "X"
$ws.Range("Q14").Value = "X"
# R14: "2" -> blank
$ws.Range("R14").ClearContents()

# --- Data validation rework --------------------------------------------------
# The old strict list rule spanned E11:R14 (including Q14). Replace it with a
# rule that covers everything except Q14, plus a dedicated rule for Q14 that
# also allows the new "X" option.
$listFormula = "1+,1,1-,2+,2,2-,3+,3,3-,4+,4,4-,5+,5,5-,6,nb,nt,t,*,/"
$listFormulaWithX = "1+,1,1-,2+,2,2-,3+,3,3-,4+,4,4-,5+,5,5-,6,nb,nt,t,*,/,X"

$ws.Range("E11:R14").Validation.Delete()

$ws.Range("E11:R13").Validation.Add(3, 1, 1, $listFormula)
$ws.Range("E11:R13").Validation.IgnoreBlank = $false
$ws.Range("E11:R13").Validation.InCellDropdown = $true
$ws.Range("E11:R13").Validation.ShowInput = $true
$ws.Range("E11:R13").Validation.ShowError = $true

$ws.Range("E14:P14").Validation.Add(3, 1, 1, $listFormula)
$ws.Range("E14:P14").Validation.IgnoreBlank = $false
$ws.Range("E14:P14").Validation.InCellDropdown = $true
$ws.Range("E14:P14").Validation.ShowInput = $true
$ws.Range("E14:P14").Validation.ShowError = $true

$ws.Range("R14").Validation.Add(3, 1, 1, $listFormula)
$ws.Range("R14").Validation.IgnoreBlank = $false
$ws.Range("R14").Validation.InCellDropdown = $true
$ws.Range("R14").Validation.ShowInput = $true
$ws.Range("R14").Validation.ShowError = $true

$ws.Range("Q14").Validation.Add(3, 1, 1, $listFormulaWithX)
$ws.Range("Q14").Validation.IgnoreBlank = $false
$ws.Range("Q14").Validation.InCellDropdown = $true
$ws.Range("Q14").Validation.ShowInput = $true
$ws.Range("Q14").Validation.ShowError = $true

# --- Selection update ---------------------------------------------------------
$ws.Range("Q15").Select()
